$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the table name value in B2 from "jenkins_template" to "jenkins_temp"
$ws.Range("B2").Value = "jenkins_temp"

# Update the active selection from C3 to B2
$ws.Range("B2").Select()
